$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp string (07:05 -> 07:35)
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 07:35"

# 2) Update Bulgaria's stats (row 82) - pure data refresh, no row shift
$ws.Range("B82").Value = 2433
$ws.Range("C82").Value = 6
$ws.Range("D82").Value = 862
$ws.Range("E82").Value = 1441

# 3) Insert a new row for "Nepal" right after "Sierra Leona" (row 125),
#    shifting "Republica del Chad" .. old "Nepal" row down by one.
$ws.Rows.Item(126).Insert()

# 4) Populate the newly inserted row with Nepal's figures
$ws.Cells.Item(126, 1).Value = "Nepal"
$ws.Cells.Item(126, 2).Value = 675
$ws.Cells.Item(126, 3).Value = 72
$ws.Cells.Item(126, 4).Value = 87
$ws.Cells.Item(126, 5).Value = 585
$ws.Cells.Item(126, 6).Value = 0
$ws.Cells.Item(126, 7).Value = 0
$ws.Cells.Item(126, 8).Value = 3

# 5) Remove the now-duplicated old "Nepal" row (originally row 131, now row 132
#    after the insert above) since its data has been relocated to row 126.
$ws.Rows.Item(132).Delete()
